# Insert a new weekly price record for Espinaca (Agrícola del Norte S.A. de
# Arica) at row 71. Excel shifts the previously-existing rows 71:124 down to
# 72:125 and extends the sheet's used range to A1:R125 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(71).Insert()

$ws.Range("A71").Value = 1
$ws.Range("B71").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C71").Value = "Arica y Parinacota"
$ws.Range("D71").Value = 45271
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = 100112012
$ws.Range("G71").Value = "Espinaca"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 400
$ws.Range("K71").Value = 800
$ws.Range("L71").Value = 1000
$ws.Range("M71").Value = 875
$ws.Range("N71").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 292
$ws.Range("Q71").Value = 3
$ws.Range("R71").Value = "Hortaliza"

# Keep the new date cell formatted the same way as the rest of column D.
$ws.Range("D71").NumberFormat = $ws.Range("D72").NumberFormat
